# Generate Report for Handback
# - Flips the two translation rows (zh-cn / de-de sheets) from
#   "Ready for handoff" to "Handed back: in sync with en-US"
# - Records the new "Latest Target File" / "Latest Handback File" hyperlink
#   columns (E/F) pointing at the same md/xlf artifacts already used for
#   handoff (A/C)
# - Stamps the "Latest Handback DateTime" column (G) with the actual
#   handback timestamp instead of the 0001-01-01 placeholder

$wb = $excel.ActiveWorkbook

$sheets = @(
    @{
        Name = "zh-cn"
        HandoffMdTarget  = "https://github.com/OpenLocalizationTest/oltest/blob/8f5de94eff47a782e6da1739044521c73a780fb1/e2e/2562a770-be01-4f71-865a-b9ee2ceafdb0.md"
        HandoffMd2Target = "https://github.com/OpenLocalizationTest/oltest/blob/8f5de94eff47a782e6da1739044521c73a780fb1/e2e/7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.md"
        XlfTarget        = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9659aee4e3bfafef5084d7b2ba16cf4e379f062d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2562a770-be01-4f71-865a-b9ee2ceafdb0.320c55bfbf76309380b4a3865348270add5aa2d6.zh-cn.xlf"
        Xlf2Target       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9659aee4e3bfafef5084d7b2ba16cf4e379f062d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.23c1d1133c193323b5962f2670ae0edb1bff762f.zh-cn.xlf"
        MdDisplay        = "2562a770-be01-4f71-865a-b9ee2ceafdb0.md"
        Md2Display       = "7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.md"
        XlfDisplay       = "2562a770-be01-4f71-865a-b9ee2ceafdb0.320c55bfbf76309380b4a3865348270add5aa2d6.zh-cn.xlf"
        Xlf2Display      = "7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.23c1d1133c193323b5962f2670ae0edb1bff762f.zh-cn.xlf"
        HandbackDateTime = "2016-03-07 02:51:33"
    },
    @{
        Name = "de-de"
        HandoffMdTarget  = "https://github.com/OpenLocalizationTest/oltest/blob/8f5de94eff47a782e6da1739044521c73a780fb1/e2e/2562a770-be01-4f71-865a-b9ee2ceafdb0.md"
        HandoffMd2Target = "https://github.com/OpenLocalizationTest/oltest/blob/8f5de94eff47a782e6da1739044521c73a780fb1/e2e/7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.md"
        XlfTarget        = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/de91d1bc00938efcc2c7b88ccddd7319aac468ea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2562a770-be01-4f71-865a-b9ee2ceafdb0.320c55bfbf76309380b4a3865348270add5aa2d6.de-de.xlf"
        Xlf2Target       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/de91d1bc00938efcc2c7b88ccddd7319aac468ea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.23c1d1133c193323b5962f2670ae0edb1bff762f.de-de.xlf"
        MdDisplay        = "2562a770-be01-4f71-865a-b9ee2ceafdb0.md"
        Md2Display       = "7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.md"
        XlfDisplay       = "2562a770-be01-4f71-865a-b9ee2ceafdb0.320c55bfbf76309380b4a3865348270add5aa2d6.de-de.xlf"
        Xlf2Display      = "7a3b1b1b-d169-4cc0-bfe1-1bf5514a98b4.23c1d1133c193323b5962f2670ae0edb1bff762f.de-de.xlf"
        HandbackDateTime = "2016-03-07 02:51:53"
    }
)

foreach ($cfg in $sheets) {
    $ws = $wb.Worksheets.Item($cfg.Name)

    # Status column (B) for the two real source-file rows: handed back.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Latest Target File (E) / Latest Handback File (F) columns now that the
    # handback round-trip produced a target + handback artifact for each row.
    $ws.Hyperlinks.Add($ws.Range("E2"), $cfg.HandoffMdTarget, "", "", $cfg.MdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F2"), $cfg.XlfTarget, "", "", $cfg.XlfDisplay)
    $ws.Hyperlinks.Add($ws.Range("E3"), $cfg.HandoffMd2Target, "", "", $cfg.Md2Display)
    $ws.Hyperlinks.Add($ws.Range("F3"), $cfg.Xlf2Target, "", "", $cfg.Xlf2Display)

    # Latest Handback DateTime (G): real timestamp instead of the
    # 0001-01-01 00:00:00 placeholder.
    $ws.Range("G2").Value = $cfg.HandbackDateTime
    $ws.Range("G3").Value = $cfg.HandbackDateTime
}
